$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data describing the "/jyt/send_email" API (row 10),
# matching the style/format already used by the other rows in the table.
$ws.Range("B10").Value = "/jyt/send_email"
$ws.Range("C10").Value = "email:String"
$ws.Range("D10").Value = "POST"

# Move the view/selection down to the newly added row, as in the saved file.
$ws.Range("D10").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
